$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 935
$ws.Range("C2").Value = 880
$ws.Range("D2").Value = 930
$ws.Range("E2").Value = 980
$ws.Range("F2").Value = 980
$ws.Range("G2").Value = 130
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 774

$ws.Range("G2").Select()
